$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.835941000000001
$ws.Range("H2").Value = 17.507823
$ws.Range("I2").Value = 0.03643643319117328
$ws.Range("J2").Value = 0.03643643319117327
$ws.Range("M2").Value = 0.008446
$ws.Range("N2").Value = 0.025338
$ws.Range("O2").Value = 0.001125187475737063
$ws.Range("P2").Value = 0.001125187475737063
$ws.Range("Q2").Value = 0.04929035768600001
$ws.Range("R2").Value = 0.443613219174
$ws.Range("S2").Value = 0.00004099781828723838
$ws.Range("T2").Value = 0.00004099781828723838
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.835941000000001
$ws.Range("H3").Value = 17.507823
$ws.Range("I3").Value = 0.03643643319117328
$ws.Range("J3").Value = 0.03643643319117327
$ws.Range("N3").Value = 9.970262999999999
$ws.Range("O3").Value = 0.4427506139949732
$ws.Range("P3").Value = 0.4427506139949733
$ws.Range("Q3").Value = 19.395288874161
$ws.Range("R3").Value = 174.557599867449
$ws.Range("S3").Value = 0.01613225316717879
$ws.Range("T3").Value = 0.01613225316717879
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.835941000000001
$ws.Range("H4").Value = 17.507823
$ws.Range("I4").Value = 0.03643643319117328
$ws.Range("J4").Value = 0.03643643319117327
$ws.Range("M4").Value = 4.174437666666667
$ws.Range("N4").Value = 12.523313
$ws.Range("O4").Value = 0.5561241985292896
$ws.Range("P4").Value = 0.5561241985292896
$ws.Range("Q4").Value = 24.36177193084434
$ws.Range("R4").Value = 219.255947377599
$ws.Range("S4").Value = 0.02026318220570724
$ws.Range("T4").Value = 0.02026318220570724
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.1093102818770573
$ws.Range("J5").Value = 0.1093102818770573
$ws.Range("M5").Value = 0.008446
$ws.Range("N5").Value = 0.025338
$ws.Range("O5").Value = 0.001125187475737063
$ws.Range("P5").Value = 0.001125187475737063
$ws.Range("Q5").Value = 0.1478724018953333
$ws.Range("R5").Value = 1.330851617058
$ws.Range("S5").Value = 0.0001229945601373529
$ws.Range("T5").Value = 0.0001229945601373529
$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.1093102818770573
$ws.Range("J6").Value = 0.1093102818770573
$ws.Range("N6").Value = 9.970262999999999
$ws.Range("O6").Value = 0.4427506139949732
$ws.Range("P6").Value = 0.4427506139949733
$ws.Range("Q6").Value = 58.18638950738701
$ws.Range("R6").Value = 523.6775055664831
$ws.Range("S6").Value = 0.04839719441703071
$ws.Range("T6").Value = 0.04839719441703071
$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.1093102818770573
$ws.Range("J7").Value = 0.1093102818770573
$ws.Range("M7").Value = 4.174437666666667
$ws.Range("N7").Value = 12.523313
$ws.Range("O7").Value = 0.5561241985292896
$ws.Range("P7").Value = 0.5561241985292896
$ws.Range("Q7").Value = 73.0859725707259
$ws.Range("R7").Value = 657.7737531365331
$ws.Range("S7").Value = 0.06079009289988922
$ws.Range("T7").Value = 0.06079009289988922
$ws.Range("G8").Value = 136.8238143333333
$ws.Range("H8").Value = 410.471443
$ws.Range("I8").Value = 0.8542532849317694
$ws.Range("J8").Value = 0.8542532849317694
$ws.Range("M8").Value = 0.008446
$ws.Range("N8").Value = 0.025338
$ws.Range("O8").Value = 0.001125187475737063
$ws.Range("P8").Value = 0.001125187475737063
$ws.Range("Q8").Value = 1.155613935859333
$ws.Range("R8").Value = 10.400525422734
$ws.Range("S8").Value = 0.0009611950973124714
$ws.Range("T8").Value = 0.0009611950973124714
$ws.Range("G9").Value = 136.8238143333333
$ws.Range("H9").Value = 410.471443
$ws.Range("I9").Value = 0.8542532849317694
$ws.Range("J9").Value = 0.8542532849317694
$ws.Range("N9").Value = 9.970262999999999
$ws.Range("O9").Value = 0.4427506139949732
$ws.Range("P9").Value = 0.4427506139949733
$ws.Range("Q9").Value = 454.723137855501
$ws.Range("R9").Value = 4092.508240699509
$ws.Range("S9").Value = 0.3782211664107638
$ws.Range("T9").Value = 0.3782211664107638
$ws.Range("G10").Value = 136.8238143333333
$ws.Range("H10").Value = 410.471443
$ws.Range("I10").Value = 0.8542532849317694
$ws.Range("J10").Value = 0.8542532849317694
$ws.Range("M10").Value = 4.174437666666667
$ws.Range("N10").Value = 12.523313
$ws.Range("O10").Value = 0.5561241985292896
$ws.Range("P10").Value = 0.5561241985292896
$ws.Range("Q10").Value = 571.1624842500732
$ws.Range("R10").Value = 5140.462358250659
$ws.Range("S10").Value = 0.4750709234236932
$ws.Range("T10").Value = 0.4750709234236932
